# chore: adapt column header formatting to respective input file names
#
# The sheet has 21 header cells in row 1 (A1:U1). Ten of them end in
# "_old" (the "FV2210" / left-hand side of the AHB diff) and ten end in
# "_new" (the "FV2304" / right-hand side); "diff" (K1) has no suffix and
# stays as-is. Rename the suffixes to the concrete format-version names,
# then wrap the sheet in an Excel Table ("Table1") and freeze the header
# row, matching the regenerated export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixMap = @{ "_old" = "_FV2210"; "_new" = "_FV2304" }

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = $cell.Value2
    foreach ($oldSuffix in $suffixMap.Keys) {
        if ($header -like "*$oldSuffix") {
            $base = $header.Substring(0, $header.Length - $oldSuffix.Length)
            $cell.Value = $base + $suffixMap[$oldSuffix]
            break
        }
    }
}

# Turn the header + data range into an Excel Table ("Table1"), matching the
# author's new xl/tables/table1.xml part (autofilter + plain table style).
$tableRange = $ws.UsedRange
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (split below row 1, bottom-left pane active), like
# the updated sheetView.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
